$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the labeling on the Shadow Wars row: "Shadow War 7PM" -> "Shadow War 6PM"
$ws.Range("A12").Value = "Shadow War 6PM"

# Reflect the user's final cursor/selection position after the edit
$ws.Range("A13").Select()
